$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-assert F1 as blank to avoid a round-trip artifact turning it into a
# non-empty value when the workbook is re-saved.
$ws.Range("F1").Value = ""

# Add "NA" values in column E (duplicate_image_filename) for rows 2 through 21
for ($r = 2; $r -le 21; $r++) {
    $ws.Range("E$r").Value = "NA"
}
